{"js": "// Update the division problems in the practice table.\n// Each entry maps the original \"A\u00f7B=\" text to its replacement \"C\u00f7D=\".\nconst replacements = [\n  [\"206\u00f77=\", \"790\u00f78=\"],\n  [\"206\u00f78=\", \"142\u00f72=\"],\n  [\"821\u00f76=\", \"180\u00f79=\"],\n  [\"340\u00f77=\", \"132\u00f72=\"],\n  [\"394\u00f76=\", \"332\u00f73=\"],\n  [\"177\u00f75=\", \"829\u00f76=\"],\n  [\"570\u00f79=\", \"378\u00f72=\"],\n  [\"265\u00f74=\", \"516\u00f73=\"],\n  [\"663\u00f72=\", \"275\u00f78=\"],\n  [\"168\u00f77=\", \"981\u00f77=\"],\n  [\"775\u00f75=\", \"225\u00f77=\"],\n  [\"706\u00f76=\", \"452\u00f76=\"],\n  [\"279\u00f76=\", \"886\u00f79=\"],\n  [\"519\u00f72=\", \"845\u00f77=\"],\n  [\"984\u00f75=\", \"167\u00f79=\"],\n  [\"120\u00f72=\", \"786\u00f76=\"],\n  [\"635\u00f76=\", \"642\u00f72=\"],\n  [\"782\u00f72=\", \"761\u00f75=\"],\n  [\"807\u00f76=\", \"613\u00f72=\"],\n  [\"929\u00f75=\", \"963\u00f72=\"],\n  [\"120\u00f76=\", \"660\u00f75=\"],\n  [\"101\u00f78=\", \"718\u00f78=\"],\n  [\"948\u00f72=\", \"687\u00f74=\"],\n  [\"933\u00f79=\", \"584\u00f75=\"],\n  [\"104\u00f76=\", \"915\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice table.\n# Each entry maps the original \"A\u00f7B=\" text to its replacement \"C\u00f7D=\".\n$replacements = @(\n    @{ Old = \"206\u00f77=\"; New = \"790\u00f78=\" },\n    @{ Old = \"206\u00f78=\"; New = \"142\u00f72=\" },\n    @{ Old = \"821\u00f76=\"; New = \"180\u00f79=\" },\n    @{ Old = \"340\u00f77=\"; New = \"132\u00f72=\" },\n    @{ Old = \"394\u00f76=\"; New = \"332\u00f73=\" },\n    @{ Old = \"177\u00f75=\"; New = \"829\u00f76=\" },\n    @{ Old = \"570\u00f79=\"; New = \"378\u00f72=\" },\n    @{ Old = \"265\u00f74=\"; New = \"516\u00f73=\" },\n    @{ Old = \"663\u00f72=\"; New = \"275\u00f78=\" },\n    @{ Old = \"168\u00f77=\"; New = \"981\u00f77=\" },\n    @{ Old = \"775\u00f75=\"; New = \"225\u00f77=\" },\n    @{ Old = \"706\u00f76=\"; New = \"452\u00f76=\" },\n    @{ Old = \"279\u00f76=\"; New = \"886\u00f79=\" },\n    @{ Old = \"519\u00f72=\"; New = \"845\u00f77=\" },\n    @{ Old = \"984\u00f75=\"; New = \"167\u00f79=\" },\n    @{ Old = \"120\u00f72=\"; New = \"786\u00f76=\" },\n    @{ Old = \"635\u00f76=\"; New = \"642\u00f72=\" },\n    @{ Old = \"782\u00f72=\"; New = \"761\u00f75=\" },\n    @{ Old = \"807\u00f76=\"; New = \"613\u00f72=\" },\n    @{ Old = \"929\u00f75=\"; New = \"963\u00f72=\" },\n    @{ Old = \"120\u00f76=\"; New = \"660\u00f75=\" },\n    @{ Old = \"101\u00f78=\"; New = \"718\u00f78=\" },\n    @{ Old = \"948\u00f72=\"; New = \"687\u00f74=\" },\n    @{ Old = \"933\u00f79=\"; New = \"584\u00f75=\" },\n    @{ Old = \"104\u00f76=\"; New = \"915\u00f79=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
